# Apply updated cryptocurrency price/volume data from the latest crawl.
# Values are written with a leading apostrophe (text qualifier) so that
# Excel keeps numeric-looking strings (e.g. "605.57", "3.227.30") as literal
# text instead of auto-converting them to numbers/dates, matching the
# original inlineStr text cells. The Style reset clears the quote-prefix
# formatting flag that gets set as a side effect, so no cell keeps a style
# index it did not have before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.737.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.31%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.237.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.86%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.44%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'158.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.234.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.77%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.86%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.89%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.76%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.505"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.52%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.45%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'39.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.91%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.763.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'66.717.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'Polkadot"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'7.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.233.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.67%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.93%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'510.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'15.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.95%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.85%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.77%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'14.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.73%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.155"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +69.56%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.92%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'28.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -6.46%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.52%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'PEPE"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.0₃0810"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +17.45%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Bittensor"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'504.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.45%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'OKB"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'55.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.79%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'dogwifhat"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0423"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.29%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.128"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.88%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.57%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -4.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Fetch.AI"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.97%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Maker"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.949.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.40%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'28.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.82%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.14%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.44%  "
$ws.Range("E51").Style = "Normal"
